# Change "project" to "task" in the "monitor actual vs planned project
# durations..." bullet on slide 3 (Content Placeholder 2), splitting the
# original single run into three runs the same way PowerPoint does when a
# word in the middle of a run is edited in place.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Locate the paragraph that contains the sentence being edited.
$targetParaIndex = 0
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $paraText = $tr.Paragraphs($i).Text
    if ($paraText -like "*monitor actual vs planned project durations*") {
        $targetParaIndex = $i
        break
    }
}

$para = $tr.Paragraphs($targetParaIndex)

# Find the "planned project " segment inside the paragraph and replace just
# the "project" word with "task", leaving the rest of the run text (and its
# formatting) untouched on either side. Editing this inner sub-range is what
# causes PowerPoint to split the single run into three runs: the unchanged
# leading text, the edited "planned task " text, and the unchanged trailing
# text.
$oldSegment = "planned project "
$newSegment = "planned task "
$segStart = $para.Text.IndexOf($oldSegment) + 1

$segRange = $para.Characters($segStart, $oldSegment.Length)
$segRange.Text = $newSegment
